$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
$ws.Range("D2").Value = "'257.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.156"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06067"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.725"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.450"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.360"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.7979"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1578"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.08046"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03345"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'0.09308"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.899"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001692"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04850"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0006161"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006207"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001100"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.003380"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001503"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.687"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.259"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3356"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1271"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003021"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04565"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").Value = "'0.009995"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.002975"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00005916"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.7514"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.06812"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00001503"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.01012"
$ws.Range("D51").Style = "Normal"

# --- Row 42/43: BKEXToken and CEJI swapped position (rank changed) ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1113"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
